$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.125.43'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '2.556.45'
$ws.Range('E3').Value = '  +3.67%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '567.99'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.33%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '146.47'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.61%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.585'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.45%  '
$ws.Range('D9').Value = '2.552.45'
$ws.Range('E9').Value = '  +3.58%  '
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.49'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.19%  '
$ws.Range('D15').Value = '3.012.30'
$ws.Range('E15').Value = '  +3.56%  '
$ws.Range('D16').Value = '63.070.98'
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000143'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.63%  '
$ws.Range('D18').Value = '2.555.51'
$ws.Range('E18').Value = '  +3.55%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.42'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.84%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '335.59'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.39%  '
$ws.Range('E21').Value = '  +1.86%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.81'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.20'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.63'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +9.39%  '
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.46'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.69%  '
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.46'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +5.65%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.32'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +7.44%  '
$ws.Range('D31').Value = '0.0₃0818'
$ws.Range('E31').Value = '  +2.96%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.85'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.56%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '175.69'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.37%  '
$ws.Range('E34').Value = '  +4.25%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '407.18'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +8.86%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.97'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('B38').Value = 'USDe'
$ws.Range('C38').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.37'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.74'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.04%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  -2.47%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '152.98'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.00%  '
$ws.Range('E44').Value = '  +2.26%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '20.99'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.21%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.606'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.73%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0526'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.83%  '
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('E49').Value = '  +5.28%  '
$ws.Range('E50').Value = '  +2.84%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.77'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.09%  '
